# Update the table style on each "Data Sources from LFX" table from the
# custom style {74CF3CD3-5AD3-4896-9624-B6345D4F0AB6} to
# {C320EC70-FFC7-4695-BD8A-B35883984C6E}.
$p = $ppt.ActivePresentation

$oldStyleId = "{74CF3CD3-5AD3-4896-9624-B6345D4F0AB6}"
$newStyleId = "{C320EC70-FFC7-4695-BD8A-B35883984C6E}"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $s = $p.Slides.Item($i)
    for ($j = 1; $j -le $s.Shapes.Count; $j++) {
        $sh = $s.Shapes.Item($j)
        if ($sh.HasTable) {
            $tbl = $sh.Table
            if ($tbl.Style -eq $oldStyleId) {
                $tbl.ApplyStyle($newStyleId)
            }
        }
    }
}
